$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.490.09"
$ws.Range("D3").Value = "1.823.28"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("D5").Value = "'317.28"
$ws.Range("E5").Value = "  +0.27%  "
$ws.Range("E6").Value = "  +0.01%  "
$ws.Range("D7").Value = "'0.5174"
$ws.Range("E7").Value = "  -2.56%  "
$ws.Range("E8").Value = "  -2.20%  "
$ws.Range("D9").Value = "'0.08513"
$ws.Range("E9").Value = "  +9.87%  "
$ws.Range("B10").Value = "Polygon"
$ws.Range("C10").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D10").Value = "'1.114"
$ws.Range("E10").Value = "  -0.28%  "
$ws.Range("B11").Value = "OKB"
$ws.Range("C11").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D11").Value = "'41.85"
$ws.Range("E11").Value = "  -0.47%  "
$ws.Range("D12").Value = "'6.446"
$ws.Range("E12").Value = "  +2.02%  "
$ws.Range("D13").Value = "'21.05"
$ws.Range("E13").Value = "  -0.43%  "
$ws.Range("E14").Value = "  -0.02%  "
$ws.Range("D15").Value = "'7.526"
$ws.Range("E15").Value = "  -0.68%  "
$ws.Range("D16").Value = "1.809.55"
$ws.Range("E16").Value = "  -1.13%  "
$ws.Range("D17").Value = "'0.00001143"
$ws.Range("E17").Value = "  +4.94%  "
$ws.Range("D18").Value = "'92.86"
$ws.Range("E18").Value = "  -0.34%  "
$ws.Range("D19").Value = "'0.06594"
$ws.Range("E19").Value = "  -0.29%  "
$ws.Range("D20").Value = "'17.78"
$ws.Range("E20").Value = "  -0.07%  "
$ws.Range("E21").Value = "  -0.04%  "
$ws.Range("D22").Value = "'6.094"
$ws.Range("E22").Value = "  +0.04%  "
$ws.Range("D23").Value = "28.520.20"
$ws.Range("E23").Value = "  -0.31%  "
$ws.Range("D24").Value = "'11.43"
$ws.Range("E24").Value = "  +2.13%  "
$ws.Range("D25").Value = "'2.272"
$ws.Range("E25").Value = "  +1.66%  "
$ws.Range("D26").Value = "'21.04"
$ws.Range("E26").Value = "  +1.22%  "
$ws.Range("D27").Value = "'159.38"
$ws.Range("E27").Value = "  +1.71%  "
$ws.Range("D28").Value = "2.025.54"
$ws.Range("E28").Value = "  -0.81%  "
$ws.Range("D29").Value = "'2.401"
$ws.Range("E29").Value = "  -0.64%  "
$ws.Range("D30").Value = "'125.51"
$ws.Range("E30").Value = "  -0.07%  "
$ws.Range("D31").Value = "'0.1091"
$ws.Range("E31").Value = "  -2.97%  "
$ws.Range("D32").Value = "'1.097"
$ws.Range("E32").Value = "  -4.56%  "
$ws.Range("D33").Value = "'5.730"
$ws.Range("E33").Value = "  -0.14%  "
$ws.Range("D34").Value = "'0.07435"
$ws.Range("E34").Value = "  +1.50%  "
$ws.Range("D35").Value = "'3.650"
$ws.Range("E35").Value = "  -0.27%  "
$ws.Range("E36").Value = "  -1.66%  "
$ws.Range("D37").Value = "'0.02358"
$ws.Range("E37").Value = "  +0.27%  "
$ws.Range("E38").Value = "  +0.27%  "
$ws.Range("D39").Value = "'8.825"
$ws.Range("E39").Value = "  -0.97%  "
$ws.Range("D40").Value = "'0.6319"
$ws.Range("E40").Value = "  +0.36%  "
$ws.Range("D41").Value = "'11.32"
$ws.Range("E41").Value = "  -0.72%  "
$ws.Range("D42").Value = "'1.195"
$ws.Range("E42").Value = "  +0.07%  "
$ws.Range("D43").Value = "'1.400"
$ws.Range("E43").Value = "  +0.10%  "
$ws.Range("D44").Value = "'13.56"
$ws.Range("E44").Value = "  -0.17%  "
$ws.Range("D45").Value = "'3.784"
$ws.Range("E45").Value = "  +1.68%  "
$ws.Range("D46").Value = "'0.5955"
$ws.Range("E46").Value = "  +0.28%  "
$ws.Range("D47").Value = "'126.42"
$ws.Range("E47").Value = "  +0.63%  "
$ws.Range("D48").Value = "'1.991"
$ws.Range("E48").Value = "  -0.40%  "
$ws.Range("D49").Value = "'1.207"
$ws.Range("E49").Value = "  +1.23%  "
$ws.Range("D50").Value = "'0.06975"
$ws.Range("E50").Value = "  +0.29%  "
$ws.Range("D51").Value = "'74.34"
$ws.Range("E51").Value = "  -0.28%  "
